$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "live-edit": update cell B2's value
$ws.Range("B2").Value = "Name + авава* dfd"

# Fix the "live-delete": update the active selection
$ws.Range("D6").Select()
